$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: "V" + "amos direto ao ponto e fazer isso funcionar " -> one run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Vamos direto ao ponto e fazer isso funcionar ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Vamos direto ao ponto e fazer isso funcionar ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Hunk 2: " 1. " + "Você precisa atuar em um projeto existente ou " -> one run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " 1. Você precisa atuar em um projeto existente ou ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " 1. Você precisa atuar em um projeto existente ou ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Hunk 3: "Clique em " (before "File > Open Library > Other") becomes
# "No FinalCutPro clique em "
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Clique em File*Open Library*Other*") {
        $r = $p.Range
        $r.Find.Execute(
            "Clique em ",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "No FinalCutPro clique em ", 2) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# Hunk 4: "[Opção " + "B" + "] " + "C" + "riar projeto..." ->
# "[Opção B] " (bold run) + "Criar projeto com resolução personalizada" (plain run)
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Opção*resolu*personalizada*") {
        $r = $p.Range
        $r.Find.Execute(
            "[Opção B] ",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "[Opção B] ", 2) | Out-Null

        $r2 = $p.Range
        $r2.Find.Execute(
            "Criar projeto com resolução personalizada",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "Criar projeto com resolução personalizada", 2) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# Hunk 5: "[" + "Command + 9" + "]" -> one run "[Command + 9]"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "[Command + 9]",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[Command + 9]", 2) | Out-Null
